# Refresh the presence (time-clock) log with the latest punches.
# The "Heure" (time) / "Employé" / "Fonction" values in each section
# (Départ / Fin pause / Pause / Arrivée) are updated in place; the
# "Type" column and the table shape are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = '17:55'
$ws.Range("B4").Value = 'Adjovi Abla'
$ws.Range("C4").Value = 'Developpeur'
$ws.Range("D4").Value = 'Départ'
$ws.Range("A5").Value = '17:54'
$ws.Range("B5").Value = 'Napo Kuvor'
$ws.Range("C5").Value = 'Mecanicien'
$ws.Range("D5").Value = 'Départ'
$ws.Range("A6").Value = '17:48'
$ws.Range("B6").Value = 'Ibrahim Kodjo'
$ws.Range("C6").Value = 'Web'
$ws.Range("D6").Value = 'Départ'
$ws.Range("A7").Value = '17:43'
$ws.Range("B7").Value = 'Amah Kwatcha'
$ws.Range("C7").Value = 'CEO'
$ws.Range("D7").Value = 'Départ'
$ws.Range("A8").Value = '16:48'
$ws.Range("B8").Value = 'Romuald Gagnon'
$ws.Range("C8").Value = 'PHOTOGRAPHE'
$ws.Range("D8").Value = 'Départ'
$ws.Range("A9").Value = '16:39'
$ws.Range("B9").Value = 'Abalo Afi'
$ws.Range("C9").Value = 'CEO'
$ws.Range("D9").Value = 'Départ'
$ws.Range("A10").Value = '16:23'
$ws.Range("B10").Value = 'Epiphanie Adoboè'
$ws.Range("C10").Value = 'CEO'
$ws.Range("D10").Value = 'Départ'
$ws.Range("A11").Value = '16:19'
$ws.Range("B11").Value = 'Akouété Kangnivi'
$ws.Range("C11").Value = 'Juriste'
$ws.Range("D11").Value = 'Départ'
$ws.Range("A12").Value = '14:54'
$ws.Range("B12").Value = 'Napo Kuvor'
$ws.Range("C12").Value = 'Mecanicien'
$ws.Range("D12").Value = 'Fin pause'
$ws.Range("A13").Value = '14:43'
$ws.Range("B13").Value = 'Amah Kwatcha'
$ws.Range("C13").Value = 'CEO'
$ws.Range("D13").Value = 'Fin pause'
$ws.Range("A14").Value = '13:07'
$ws.Range("B14").Value = 'Abi Conrad'
$ws.Range("C14").Value = 'CEO'
$ws.Range("D14").Value = 'Fin pause'
$ws.Range("A15").Value = '12:54'
$ws.Range("B15").Value = 'Napo Kuvor'
$ws.Range("C15").Value = 'Mecanicien'
$ws.Range("D15").Value = 'Pause'
$ws.Range("A16").Value = '12:43'
$ws.Range("B16").Value = 'Amah Kwatcha'
$ws.Range("C16").Value = 'CEO'
$ws.Range("D16").Value = 'Pause'
$ws.Range("A17").Value = '11:07'
$ws.Range("B17").Value = 'Abi Conrad'
$ws.Range("C17").Value = 'CEO'
$ws.Range("D17").Value = 'Pause'
$ws.Range("A18").Value = '08:55'
$ws.Range("B18").Value = 'Adjovi Abla'
$ws.Range("C18").Value = 'Developpeur'
$ws.Range("D18").Value = 'Arrivée'
$ws.Range("A19").Value = '08:54'
$ws.Range("B19").Value = 'Napo Kuvor'
$ws.Range("C19").Value = 'Mecanicien'
$ws.Range("D19").Value = 'Arrivée'
$ws.Range("A20").Value = '08:48'
$ws.Range("B20").Value = 'Ibrahim Kodjo'
$ws.Range("C20").Value = 'Web'
$ws.Range("D20").Value = 'Arrivée'
$ws.Range("A21").Value = '08:43'
$ws.Range("B21").Value = 'Amah Kwatcha'
$ws.Range("C21").Value = 'CEO'
$ws.Range("D21").Value = 'Arrivée'
$ws.Range("A22").Value = '07:48'
$ws.Range("B22").Value = 'Romuald Gagnon'
$ws.Range("C22").Value = 'PHOTOGRAPHE'
$ws.Range("D22").Value = 'Arrivée'
$ws.Range("A23").Value = '07:39'
$ws.Range("B23").Value = 'Abalo Afi'
$ws.Range("C23").Value = 'CEO'
$ws.Range("D23").Value = 'Arrivée'
$ws.Range("A24").Value = '07:23'
$ws.Range("B24").Value = 'Epiphanie Adoboè'
$ws.Range("C24").Value = 'CEO'
$ws.Range("D24").Value = 'Arrivée'
$ws.Range("A25").Value = '07:19'
$ws.Range("B25").Value = 'Akouété Kangnivi'
$ws.Range("C25").Value = 'Juriste'
$ws.Range("D25").Value = 'Arrivée'
$ws.Range("A26").Value = '07:07'
$ws.Range("B26").Value = 'Abi Conrad'
$ws.Range("C26").Value = 'CEO'
$ws.Range("D26").Value = 'Arrivée'
